# Apply the new table style ("Light Style 2 - Accent 1") to the three
# tables that still carried the deck's original custom default table
# style. This mirrors picking a style from the Table Design > Table
# Styles gallery for each table in PowerPoint.

$p = $ppt.ActivePresentation

$newStyleId = "{07E38A7B-6411-43C0-B6AB-9AB1E51C437C}"
$targetSlides = @(14, 15, 16)

foreach ($slideIndex in $targetSlides) {
    $slide = $p.Slides.Item($slideIndex)
    $shape = $slide.Shapes.Item(1)
    if ($shape.HasTable) {
        $table = $shape.Table
        $table.ApplyStyle($newStyleId)
    }
}
